# Resort the sheets: "总计" (the summary sheet) should come before "2021-Q2".
# The underlying cell data/formatting of each sheet is unchanged - only the
# tab order (and therefore which sheet is first/active) is updated.
$wb = $excel.ActiveWorkbook

$summarySheet = $wb.Worksheets.Item("总计")
$quarterSheet = $wb.Worksheets.Item("2021-Q2")

# Move "总计" so it sits immediately before "2021-Q2", making it the first tab.
$summarySheet.Move($quarterSheet)
